$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("errorDefs")

# Header: F1 currently "list#anotherList[0].ccc" -> change to "list#anotherList[1].ccc"
$ws.Range("F1").Value = "list#anotherList[1].ccc"

# F5 currently "ccc2-3" -> change to "ccc3-2"
$ws.Range("F5").Value = "ccc3-2"

# Move selection to F6 to match saved workbook state
$ws.Range("F6").Select()
